$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Project box"
$ws.Range("B11").Value = "https://www.amazon.com/gp/product/B07D23BF7Y"
$ws.Range("C11").Value = 12.59

$ws.Range("C11").Select()
